$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidential/date notice text in A16
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.02803438230997498
$ws.Range("E2").Value = 0.003819223424570195

$ws.Range("D3").Value = 0.02236989222242707
$ws.Range("E3").Value = 0.0009324009324009896

$ws.Range("D4").Value = 0.05869046268538592
$ws.Range("E4").Value = -0.001155001155001223

$ws.Range("D5").Value = 0.1385716618135801
$ws.Range("E5").Value = 0.0003177629488400058

$ws.Range("D6").Value = 0.02176849496454364
$ws.Range("E6").Value = -0.005280528052805322

$ws.Range("D7").Value = 0.1266320142879937
$ws.Range("E7").Value = -0.001262785705265723

$ws.Range("D8").Value = 0.09130670414256105
$ws.Range("E8").Value = -0.004751461988304118

$ws.Range("D9").Value = 0.03145249720651553
$ws.Range("E9").Value = -0.01382298328863218

$ws.Range("D10").Value = 0.1080304147103392
$ws.Range("E10").Value = -0.009518477043673035

$ws.Range("D11").Value = 0.2840393455175568
$ws.Range("E11").Value = 0.01124291784702547

$ws.Range("D12").Value = 0.08910413013912206
$ws.Range("E12").Value = 0.006179775280898969

$ws.Range("E13").Value = 0.001676496247987824

# Restore worksheet protection (content was protected before this edit)
$ws.Protect()
